$wb = $excel.ActiveWorkbook

# This script applies the 2022-12-07 daily crime-data refresh to the
# Chicago Transit Authority violent-crime YTD workbook. One new incident
# record was added for each of 14 community areas; the new counts cascade
# into: that area's own worksheet (category row + Total row), the
# "By Neighborhood" rollup sheet, and the "Citywide Totals" sheet.

# --- Citywide Totals ---
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("B2").Value = 43
$ws.Range("G2").Value = 84
$ws.Range("E2").Value = 70
$ws.Range("F2").Value = 93
$ws.Range("C3").Value = 77
$ws.Range("I3").Value = 194
$ws.Range("H3").Value = 155
$ws.Range("B6").Value = 376
$ws.Range("H6").Value = 441
$ws.Range("F6").Value = 540
$ws.Range("I6").Value = 499
$ws.Range("G7").Value = 663
$ws.Range("E7").Value = 697
$ws.Range("F7").Value = 780
$ws.Range("C7").Value = 631
$ws.Range("B7").Value = 505
$ws.Range("I7").Value = 833
$ws.Range("H7").Value = 717

# --- By Neighborhood ---
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("B7").Value = 21
$ws.Range("H14").Value = 2
$ws.Range("F21").Value = 13
$ws.Range("F30").Value = 5
$ws.Range("F35").Value = 8
$ws.Range("I36").Value = 36
$ws.Range("G41").Value = 11
$ws.Range("E45").Value = 2
$ws.Range("B50").Value = 9
$ws.Range("H53").Value = 97
$ws.Range("C77").Value = 24
$ws.Range("I88").Value = 11
$ws.Range("F91").Value = 11
$ws.Range("I95").Value = 5
$ws.Range("E98").Value = 697
$ws.Range("B98").Value = 505
$ws.Range("I98").Value = 833
$ws.Range("H98").Value = 717
$ws.Range("G98").Value = 663
$ws.Range("F98").Value = 780
$ws.Range("C98").Value = 631

# --- Roseland ---
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("C3").Value = 4
$ws.Range("C7").Value = 24

# --- Auburn Gresham ---
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("B2").Value = 2
$ws.Range("B6").Value = 21

# --- Gage Park ---
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("F2").Value = 1
$ws.Range("F6").Value = 5

# --- Chinatown ---
$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("F6").Value = 6
$ws.Range("F7").Value = 13

# --- Grand Crossing ---
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 36

# --- Little Italy, UIC ---
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 9

# --- Washington Park ---
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I3").Value = 3
$ws.Range("I5").Value = 11

# --- Humboldt Park ---
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("G2").Value = 4
$ws.Range("G5").Value = 11

# --- Loop ---
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("H6").Value = 62
$ws.Range("H7").Value = 97

# --- West Loop ---
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("F6").Value = 9
$ws.Range("F7").Value = 11

# --- Grand Boulevard ---
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = 8

# --- Jefferson Park ---
$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("E2").Value = 1
$ws.Range("E6").Value = 2

# --- Bridgeport ---
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("G3").Value = 1
$ws.Range("G5").Value = 2

# --- Wicker Park ---
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 5
